# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Santa Lucia" / "Suazilandia" order (rows 184-185, country column A) ---
$ws.Range("A184").Value = "Suazilandia"
$ws.Range("A185").Value = "Santa Lucia"

# --- Update "Datos actualizados..." timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 21:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 605354
$ws.Range("C4").Value = 18413
$ws.Range("D4").Value = 38166
$ws.Range("E4").Value = 541794
$ws.Range("G4").Value = 1754
$ws.Range("H4").Value = 25394

# --- Row 63: Barein ---
$ws.Range("B63").Value = 1528
$ws.Range("C63").Value = 167
$ws.Range("E63").Value = 876

# --- Row 72: Uzbekistan ---
$ws.Range("B72").Value = 1165
$ws.Range("C72").Value = 167
$ws.Range("E72").Value = 1062

# --- Row 184: now Suazilandia (active/recovered updated) ---
$ws.Range("D184").Value = 8
$ws.Range("E184").Value = 7

# --- Row 185: now Santa Lucia (active/recovered updated) ---
$ws.Range("D185").Value = 11
$ws.Range("E185").Value = 4
